$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.244007110595703
$ws.Range("B1").Value = 2.331911087036133
$ws.Range("C1").Value = 3.119572877883911
$ws.Range("D1").Value = 3.599503755569458
$ws.Range("E1").Value = 1.31069540977478
